$wb = $excel.ActiveWorkbook

# Sheet: 展览 (37 cell updates to column F, "想去人数")
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 42231
$ws.Range("F3").Value = 26
$ws.Range("F4").Value = 9890
$ws.Range("F5").Value = 218
$ws.Range("F6").Value = 1005
$ws.Range("F7").Value = 955
$ws.Range("F8").Value = 766
$ws.Range("F9").Value = 235
$ws.Range("F10").Value = 315
$ws.Range("F11").Value = 1001
$ws.Range("F14").Value = 785
$ws.Range("F15").Value = 341
$ws.Range("F16").Value = 1570
$ws.Range("F18").Value = 769
$ws.Range("F19").Value = 745
$ws.Range("F20").Value = 489
$ws.Range("F21").Value = 716
$ws.Range("F22").Value = 785
$ws.Range("F24").Value = 255
$ws.Range("F25").Value = 70
$ws.Range("F26").Value = 563
$ws.Range("F27").Value = 554
$ws.Range("F28").Value = 70
$ws.Range("F29").Value = 269
$ws.Range("F30").Value = 961
$ws.Range("F32").Value = 449
$ws.Range("F33").Value = 117
$ws.Range("F34").Value = 228
$ws.Range("F35").Value = 171
$ws.Range("F36").Value = 462
$ws.Range("F37").Value = 1390
$ws.Range("F38").Value = 318
$ws.Range("F39").Value = 1300
$ws.Range("F40").Value = 386
$ws.Range("F41").Value = 104
$ws.Range("F43").Value = 49
$ws.Range("F46").Value = 5

# Sheet: 演出 (10 cell updates to column F, "想去人数")
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 221
$ws.Range("F5").Value = 4469
$ws.Range("F7").Value = 348
$ws.Range("F11").Value = 141
$ws.Range("F12").Value = 11
$ws.Range("F13").Value = 63
$ws.Range("F14").Value = 2
$ws.Range("F19").Value = 17
$ws.Range("F20").Value = 4390
$ws.Range("F23").Value = 7

# Sheet: 本地生活 (3 cell updates to column F, "想去人数")
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2084
$ws.Range("F3").Value = 559
$ws.Range("F4").Value = 459

# Sheet: 全部类型 (35 cell updates to column F, "想去人数")
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2084
$ws.Range("F3").Value = 559
$ws.Range("F5").Value = 221
$ws.Range("F6").Value = 348
$ws.Range("F7").Value = 26
$ws.Range("F9").Value = 9890
$ws.Range("F10").Value = 218
$ws.Range("F11").Value = 1005
$ws.Range("F12").Value = 1005
$ws.Range("F14").Value = 955
$ws.Range("F15").Value = 141
$ws.Range("F16").Value = 235
$ws.Range("F17").Value = 315
$ws.Range("F18").Value = 1001
$ws.Range("F19").Value = 11
$ws.Range("F21").Value = 786
$ws.Range("F22").Value = 341
$ws.Range("F23").Value = 1570
$ws.Range("F25").Value = 745
$ws.Range("F26").Value = 489
$ws.Range("F27").Value = 716
$ws.Range("F28").Value = 785
$ws.Range("F30").Value = 70
$ws.Range("F31").Value = 563
$ws.Range("F32").Value = 2
$ws.Range("F34").Value = 554
$ws.Range("F35").Value = 70
$ws.Range("F36").Value = 269
$ws.Range("F39").Value = 449
$ws.Range("F40").Value = 117
$ws.Range("F41").Value = 228
$ws.Range("F42").Value = 1300
$ws.Range("F43").Value = 386
$ws.Range("F44").Value = 104
$ws.Range("F47").Value = 49

